# Cambiando lógica a formato múltiple
# Rewrites the Huff-model result table (Sheet1) to the new "multiple format"
# layout: 5 zones x 4 stores = 20 data rows (rows 2..21), all sharing a
# common "Total SV" (column D) of 25000, and recalculated Factor / Dis (Km) /
# Captura Final / % Captura Final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data grid: row -> (Zone, Store, C, D, E, F, G, H, I)
$data = @(
    @(2,  "Z408",   "Supermercado",            2000,  25000, 0.08,                3.2, 0.26400000000000001, 9702.98,             0.66227741174109767),
    @(3,  "Z408",   "Plaza Vea Izaguirre",     10000,  25000, 0.4,                 2,   0.50900000000000001, 3798.9,              0.25929411989546058),
    @(4,  "Z408",   "Plaza Vea Los Olivos",     6000,  25000, 0.24,                1.2, 1.603,                706.93,              0.048251544441206107),
    @(5,  "Z408",   "Plaza Vea Universitaria",  7000,  25000, 0.28000000000000003,1.6, 3.48,                 442.12,              0.03017692392223565),

    @(6,  "Z407",   "Supermercado",             2000,  25000, 0.08,                3.2, 0.57099999999999995, 16816.150000000001, 0.54986202889300384),
    @(7,  "Z407",   "Plaza Vea Izaguirre",     10000,  25000, 0.4,                 2,   0.90100000000000002, 8029.09,             0.26253879262283752),
    @(8,  "Z407",   "Plaza Vea Los Olivos",     6000,  25000, 0.24,                1.2, 1.226,                3896.83,             0.12742029834719151),
    @(9,  "Z407",   "Plaza Vea Universitaria",  7000,  25000, 0.28000000000000003,1.6, 3.1789999999999998,  1840.42,             0.060178880136967268),

    @(10, "Z651",   "Supermercado",             2000,  25000, 0.08,                3.2, 0.57799999999999996, 615.92999999999995, 0.66824706252508925),
    @(11, "Z651",   "Plaza Vea Izaguirre",     10000,  25000, 0.4,                 2,   1.1060000000000001,  202.91,              0.22014516496511921),
    @(12, "Z651",   "Plaza Vea Los Olivos",     6000,  25000, 0.24,                1.2, 2.0249999999999999,  61.33,               0.066539367046034012),
    @(13, "Z651",   "Plaza Vea Universitaria",  7000,  25000, 0.28000000000000003,1.6, 3.968,                41.54,               0.045068405463757587),

    @(14, "Z414",   "Supermercado",             2000,  25000, 0.08,                3.2, 1.0960000000000001,  3744.17,             0.31923390751333702),
    @(15, "Z414",   "Plaza Vea Los Olivos",     6000,  25000, 0.24,                1.2, 0.56000000000000005, 3346.19,             0.28530149779044572),
    @(16, "Z414",   "Plaza Vea Izaguirre",     10000,  25000, 0.4,                 2,   0.68600000000000005, 3758.41,             0.32044803263131783),
    @(17, "Z414",   "Plaza Vea Universitaria",  7000,  25000, 0.28000000000000003,1.6, 2.327,                879.84,              0.075016562064899414),

    @(18, 'Z396"',  "Supermercado",             2000,  25000, 0.08,                0.8, 4.0750000000000002,  6.21,                0.47695852534562211),
    @(19, 'Z396"',  "Plaza Vea Izaguirre",     10000,  25000, 0.4,                 0.5, 4.4249999999999998,  3.54,                0.27188940092165897),
    @(20, 'Z396"',  "Plaza Vea Los Olivos",     6000,  25000, 0.24,                0.3, 5.5460000000000003,  1.66,                0.12749615975422429),
    @(21, 'Z396"',  "Plaza Vea Universitaria",  7000,  25000, 0.28000000000000003,0.4, 7.4260000000000002,  1.61,                0.1236559139784946)
)

# Write column A (Zone) for every row first, then column B (Store), so that
# newly-introduced labels are registered in the same order a manual
# zone-by-zone, then store-by-store fill would produce (Z396" before
# "Plaza Vea Universitaria").
foreach ($row in $data) {
    $r    = $row[0]
    $zone = $row[1]
    $ws.Cells.Item($r, 1).Value = $zone
}
foreach ($row in $data) {
    $r     = $row[0]
    $store = $row[2]
    $ws.Cells.Item($r, 2).Value = $store
}

foreach ($row in $data) {
    $r      = $row[0]
    $cVal   = $row[3]
    $dVal   = $row[4]
    $eVal   = $row[5]
    $fVal   = $row[6]
    $gVal   = $row[7]
    $hVal   = $row[8]
    $iVal   = $row[9]

    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
    $ws.Cells.Item($r, 6).Value = $fVal
    $ws.Cells.Item($r, 7).Value = $gVal
    $ws.Cells.Item($r, 8).Value = $hVal
    $ws.Cells.Item($r, 9).Value = $iVal
}

# Update the visible selection to match the new last block (rows 18-21),
# mirroring the original workbook's pattern of selecting the last zone block.
$ws.Range("A18:I21").Select()
